# "Corrected diagrams for undo/redo feature implementation."
#
# The slide has two textboxes whose label text ("currentStatePointer = 1"
# and "currentStatePointer = 2") was previously split across two runs:
#   run 1: "currentStatePointer"   (flagged err="1" by the spell checker)
#   run 2: " = 1" / " = 2"         (not flagged)
# The fix merges each label back into a single run - the surviving run's
# formatting is that of the second (non-flagged) run, so the err="1"
# marker goes away.
#
# Note: assigning TextRange.Text is a no-op when the new string already
# equals the concatenation of the existing runs, so instead we grow the
# *second* run so it already contains the full combined label, then
# delete the now-redundant leading characters that used to be the first
# run. That collapses the paragraph to a single run with the right text
# and formatting.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# TextBox 23: "currentStatePointer" + " = 1" -> "currentStatePointer = 1"
$sh1 = $s.Shapes.Item("TextBox 23")
$tr1 = $sh1.TextFrame.TextRange
$tr1.Characters(20, $tr1.Text.Length - 19).Text = "currentStatePointer = 1"
$sh1.TextFrame.TextRange.Characters(1, 19).Text = ""

# TextBox 29: "currentStatePointer" + " = 2" -> "currentStatePointer = 2"
$sh2 = $s.Shapes.Item("TextBox 29")
$tr2 = $sh2.TextFrame.TextRange
$tr2.Characters(20, $tr2.Text.Length - 19).Text = "currentStatePointer = 2"
$sh2.TextFrame.TextRange.Characters(1, 19).Text = ""
